# Apply the "oop in the player & new inputs :)" edit:
# insert six new song rows into the existing song table, pushing the
# trailing rows down, while leaving already-present rows' data intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert new row before "alarm.mp3" (currently row 3) ---
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "perfect-beauty-191271.mp3"
$ws.Range("B3").Value = 7.3412
$ws.Range("C3").Value = 440.472
$ws.Range("D3").Value = 440.472
$ws.Range("E3").Value = "[161.49902344]"

# --- insert new row before "testfile3.mp3" (now row 5) ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "endless-horizons-223627.mp3"
$ws.Range("B5").Value = 1.7812
$ws.Range("C5").Value = 106.872
$ws.Range("D5").Value = 106.872
$ws.Range("E5").Value = "[117.45383523]"

# --- insert new row before "testfile2.mp3" (now row 7) ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "lotus-sky-dreams-216049.mp3"
$ws.Range("B7").Value = 2.4112
$ws.Range("C7").Value = 144.672
$ws.Range("D7").Value = 144.672
$ws.Range("E7").Value = "[99.38401442]"

# --- insert three new rows before "alarm2.mp3" (now row 11) ---
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "night-detective-226857.mp3"
$ws.Range("B11").Value = 1.9332
$ws.Range("C11").Value = 115.992
$ws.Range("D11").Value = 115.992
$ws.Range("E11").Value = "[73.828125]"

$ws.Range("A12").Value = "titanium.mp3"
$ws.Range("B12").Value = 1.77
$ws.Range("C12").Value = 106.2
$ws.Range("D12").Value = 106.2
$ws.Range("E12").Value = "[64.59960938]"

$ws.Range("A13").Value = "island-breeze-214305.mp3"
$ws.Range("B13").Value = 1.5108
$ws.Range("C13").Value = 90.648
$ws.Range("D13").Value = 90.648
$ws.Range("E13").Value = "[107.66601562]"
